# Runtime_Estimates.xlsx edit
# Commit: "Added results for 10,000 and 50,000 sample size runs.
#          Updated code to print end time to file as well as console."
#
# 1) New measured data point: the 10,000-permutation run actually took
#    1360 sec (B7) once it finished (previously only the estimate existed).
# 2) With the new measured point, the estimate curve (column C) was
#    refit: every row from 3000 permutations onward (rows 5-15) now uses
#    a new polynomial instead of the old cubic.
# 3) A new scratch table was added at A22:C27 re-deriving a simple
#    permutations/runtime ratio (A/B) for the original sample sizes plus
#    the new 10,000 one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) record the measured runtime for the 10,000-permutation run ---
$ws.Range("B7").Value = 1360

# --- 2) refit the estimate polynomial in column C, rows 5 through 15 ---
$ws.Range("C5:C15").Formula = "=(-5*10^(-10))*(A5^3)+(2*10^(-5))*(A5^2)-(0.0023)*(A5) + 2.8093"

# --- 3) new scratch table: permutations / runtime ratio ---
$ws.Range("A22").Value = 100
$ws.Range("B22").Value = 2
$ws.Range("C22").Formula = "=A22/B22"

$ws.Range("A23").Value = 500
$ws.Range("B23").Value = 8
$ws.Range("A24").Value = 1000
$ws.Range("B24").Value = 18
$ws.Range("A25").Value = 3000
$ws.Range("B25").Value = 153
$ws.Range("A26").Value = 5000
$ws.Range("A27").Value = 10000
$ws.Range("B27").Value = 1360

$ws.Range("C23:C27").Formula = "=A23/B23"
# Row 26 (5000 permutations) still has no recorded runtime, so it has no
# ratio either - clear the div/0 that the fill produced for it.
$ws.Range("C26").Clear()

# Match the author's final selection (cell C27) as recorded in the diff.
$ws.Range("C27").Select()
